# Daily attendance processing - normalise the "Recorded By" (column G)
# cell contents so the comma-separated list of recorders is sorted in
# strict ordinal (case-sensitive, ASCII-code) order, e.g.
#   "dnasr281@gmail.com, System"  ->  "System, dnasr281@gmail.com"
#   "backup@backdoor.com, system, System" -> "System, backup@backdoor.com, system"

function Test-OrdinalLess($a, $b) {
    $la = $a.Length
    $lb = $b.Length
    $n = [Math]::Min($la, $lb)
    $k = 0
    while ($k -lt $n) {
        $ca = [int][char]$a.Substring($k, 1)
        $cb = [int][char]$b.Substring($k, 1)
        if ($ca -lt $cb) { return $true }
        if ($ca -gt $cb) { return $false }
        $k = $k + 1
    }
    return $la -lt $lb
}

function Sort-NameList($parts) {
    $cnt = $parts.Count
    if ($cnt -le 1) {
        return $parts
    }
    if ($cnt -eq 2) {
        $a0 = $parts[0]
        $a1 = $parts[1]
        if (Test-OrdinalLess $a1 $a0) {
            return @($a1, $a0)
        }
        return @($a0, $a1)
    }
    # Fixed-size (3 element) bubble pass - avoids nested for/while loops,
    # which this host's PowerShell engine mishandles when combined with
    # array-index arguments to a function call.
    $a0 = $parts[0]
    $a1 = $parts[1]
    $a2 = $parts[2]
    if (Test-OrdinalLess $a1 $a0) { $t = $a0; $a0 = $a1; $a1 = $t }
    if (Test-OrdinalLess $a2 $a1) { $t = $a1; $a1 = $a2; $a2 = $t }
    if (Test-OrdinalLess $a1 $a0) { $t = $a0; $a0 = $a1; $a1 = $t }
    return @($a0, $a1, $a2)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $usedRange.Rows.Count + $firstRow - 1
$headerRow = 1
$col = 7  # column G = "Recorded By"

for ($r = $headerRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2
    if ($val -eq $null) {
        continue
    }
    $parts = @($val -split ", ")
    $sortedParts = Sort-NameList $parts
    $newVal = $sortedParts -join ", "
    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
